$d = $word.ActiveDocument

$replacements = @(
    @{Old = "2024-10-21 Monday"; New = "2024-10-22 Tuesday"},
    @{Old = "78×45=3510"; New = "25×77=1925"},
    @{Old = "39×61=2379"; New = "91×68=6188"},
    @{Old = "27×92=2484"; New = "42×76=3192"},
    @{Old = "36×45=1620"; New = "44×36=1584"},
    @{Old = "14×41=574"; New = "16×26=416"},
    @{Old = "66×68=4488"; New = "87×61=5307"},
    @{Old = "89×36=3204"; New = "48×82=3936"},
    @{Old = "54×82=4428"; New = "40×39=1560"},
    @{Old = "26×70=1820"; New = "91×99=9009"},
    @{Old = "70×51=3570"; New = "45×87=3915"},
    @{Old = "53×32=1696"; New = "95×47=4465"},
    @{Old = "67×81=5427"; New = "20×50=1000"},
    @{Old = "75×65=4875"; New = "87×75=6525"},
    @{Old = "38×73=2774"; New = "47×27=1269"},
    @{Old = "47×57=2679"; New = "94×90=8460"},
    @{Old = "33×69=2277"; New = "57×24=1368"},
    @{Old = "67×28=1876"; New = "26×30=780"},
    @{Old = "97×96=9312"; New = "84×28=2352"},
    @{Old = "35×40=1400"; New = "22×89=1958"},
    @{Old = "72×17=1224"; New = "88×27=2376"},
    @{Old = "67×78=5226"; New = "40×66=2640"},
    @{Old = "96×86=8256"; New = "64×36=2304"},
    @{Old = "83×56=4648"; New = "50×42=2100"},
    @{Old = "77×31=2387"; New = "27×19=513"},
    @{Old = "81×35=2835"; New = "87×58=5046"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
